$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.414.27'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '1.817.68'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.34'
$ws.Range('E5').Value = '  -0.58%  '
$ws.Range('E6').Value = '  +0.10%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5127'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3948'
$ws.Range('E8').Value = '  -2.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07985'
$ws.Range('E9').Value = '  +4.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '41.70'
$ws.Range('E10').Value = '  -0.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.108'
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '20.96'
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.252'
$ws.Range('E13').Value = '  -1.13%  '
$ws.Range('E14').Value = '  +0.05%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.496'
$ws.Range('E15').Value = '  -1.57%  '
$ws.Range('D16').Value = '1.829.43'
$ws.Range('E16').Value = '  +0.29%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001134'
$ws.Range('E17').Value = '  +5.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '92.63'
$ws.Range('E18').Value = '  +3.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06630'
$ws.Range('E19').Value = '  +0.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.65'
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.088'
$ws.Range('E22').Value = '  -0.37%  '
$ws.Range('D23').Value = '28.444.43'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.25'
$ws.Range('E24').Value = '  +0.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.271'
$ws.Range('E25').Value = '  +3.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '21.09'
$ws.Range('E26').Value = '  +2.63%  '
$ws.Range('D27').Value = '2.036.66'
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '155.47'
$ws.Range('E28').Value = '  -1.43%  '
$ws.Range('E29').Value = '  -2.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.69'
$ws.Range('E30').Value = '  +1.61%  '
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('E32').Value = '  -1.93%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.711'
$ws.Range('E33').Value = '  +1.05%  '
$ws.Range('E34').Value = '  +0.19%  '
$ws.Range('E35').Value = '  -3.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.2222'
$ws.Range('E36').Value = '  -0.49%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02325'
$ws.Range('E37').Value = '  -0.75%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.191'
$ws.Range('E38').Value = '  +0.04%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.829'
$ws.Range('E39').Value = '  -0.41%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6252'
$ws.Range('E40').Value = '  +0.09%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '11.28'
$ws.Range('E41').Value = '  +0.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.176'
$ws.Range('E42').Value = '  -0.40%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.400'
$ws.Range('E44').Value = '  -0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.43'
$ws.Range('E45').Value = '  -0.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.740'
$ws.Range('E46').Value = '  +0.94%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5904'
$ws.Range('E47').Value = '  +1.15%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '124.79'
$ws.Range('E48').Value = '  -0.32%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.971'
$ws.Range('E49').Value = '  -0.83%  '
$ws.Range('E50').Value = '  -1.28%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06884'
